$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Professional summary: append a sentence to the end of the existing
#    italic "...teammates. Proven multitasking ability." run.
#    We locate the run with Find (without replacing) and then use
#    InsertAfter on the collapsed end of the match so the existing runs in
#    this multi-run paragraph are left completely untouched, and give the
#    newly appended text the same italic formatting as its neighbour.
# ---------------------------------------------------------------------------
$summaryRange = $d.Content
$foundSummary = $summaryRange.Find.Execute(
    "teammates. Proven multitasking ability.", $true, $false, $false,
    $false, $false, $true, 1, $false, "", 0)
if ($foundSummary) {
    $summaryRange.Collapse(0)
    $summaryRange.InsertAfter(" I have accuracy, consistency, punctuality skills.")
    $summaryRange.Font.Italic = -1
}

# ---------------------------------------------------------------------------
# 2-6) Simple one-for-one wording tweaks inside single-run bullet paragraphs.
# ---------------------------------------------------------------------------
function Replace-Text($find, $replace) {
    $range = $d.Content
    $range.Find.Execute($find, $true, $false, $false, $false, $false, $true,
                         1, $false, $replace, 2)
}

Replace-Text "Identifying bottlenecks and bugs, and devise solutions to these problems." "identity bottlenecks and bugs, and devise solutions to these problems."
Replace-Text "Developed features using Ruby on Rails, HTML, CSS and JavaScript." "Developing features using Ruby on Rails, HTML, CSS and JavaScript."
Replace-Text "Supported product migration and platform upgrades." "Support product migration and platform upgrades."
Replace-Text "Worked on implementing Bug Reporting Feature to allow customers to create reports about bugs directly to Redmine APIs." "Working on implementing Bug Reporting Feature to allow customers to create reports about bugs directly to Redmine APIs."
Replace-Text "Experience with multiple 3rd party service integrations including Braintree, Google maps." "Experienced with multiple 3rd party service integrations including Braintree, Google maps."
